$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 12.97098566666667
$ws.Cells.Item(2, 8).Value = 38.91295700000001
$ws.Cells.Item(2, 9).Value = 0.7291028508134716
$ws.Cells.Item(2, 10).Value = 0.7291028508134717
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 30.020559
$ws.Cells.Item(2, 14).Value = 90.061677
$ws.Cells.Item(2, 15).Value = 0.8829766276144534
$ws.Cells.Item(2, 16).Value = 0.8829766276144534
$ws.Cells.Item(2, 17).Value = 389.3962404943211
$ws.Cells.Item(2, 18).Value = 3504.56616444889
$ws.Cells.Item(2, 19).Value = 0.643780776395363
$ws.Cells.Item(2, 20).Value = 0.6437807763953631

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 12.97098566666667
$ws.Cells.Item(3, 8).Value = 38.91295700000001
$ws.Cells.Item(3, 9).Value = 0.7291028508134716
$ws.Cells.Item(3, 10).Value = 0.7291028508134717
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 1.426936666666667
$ws.Cells.Item(3, 14).Value = 4.28081
$ws.Cells.Item(3, 15).Value = 0.04196962907162197
$ws.Cells.Item(3, 16).Value = 0.04196962907162197
$ws.Cells.Item(3, 17).Value = 18.50877505057445
$ws.Cells.Item(3, 18).Value = 166.57897545517
$ws.Cells.Item(3, 19).Value = 0.03060017620370353
$ws.Cells.Item(3, 20).Value = 0.03060017620370354

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 12.97098566666667
$ws.Cells.Item(4, 8).Value = 38.91295700000001
$ws.Cells.Item(4, 9).Value = 0.7291028508134716
$ws.Cells.Item(4, 10).Value = 0.7291028508134717
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 1.067853
$ws.Cells.Item(4, 14).Value = 3.203559
$ws.Cells.Item(4, 15).Value = 0.03140811737476231
$ws.Cells.Item(4, 16).Value = 0.0314081173747623
$ws.Cells.Item(4, 17).Value = 13.85110595710701
$ws.Cells.Item(4, 18).Value = 124.659953613963
$ws.Cells.Item(4, 19).Value = 0.02289974791662333
$ws.Cells.Item(4, 20).Value = 0.02289974791662333

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 12.97098566666667
$ws.Cells.Item(5, 8).Value = 38.91295700000001
$ws.Cells.Item(5, 9).Value = 0.7291028508134716
$ws.Cells.Item(5, 10).Value = 0.7291028508134717
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 1.483919333333333
$ws.Cells.Item(5, 14).Value = 4.451758
$ws.Cells.Item(5, 15).Value = 0.04364562593916237
$ws.Cells.Item(5, 16).Value = 0.04364562593916237
$ws.Cells.Item(5, 17).Value = 19.24789640315623
$ws.Cells.Item(5, 18).Value = 173.231067628406
$ws.Cells.Item(5, 19).Value = 0.03182215029778169
$ws.Cells.Item(5, 20).Value = 0.0318221502977817

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 1.047813333333333
$ws.Cells.Item(6, 8).Value = 3.14344
$ws.Cells.Item(6, 9).Value = 0.05889789011308234
$ws.Cells.Item(6, 10).Value = 0.05889789011308236
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 30.020559
$ws.Cells.Item(6, 14).Value = 90.061677
$ws.Cells.Item(6, 15).Value = 0.8829766276144534
$ws.Cells.Item(6, 16).Value = 0.8829766276144534
$ws.Cells.Item(6, 17).Value = 31.45594199432
$ws.Cells.Item(6, 18).Value = 283.10347794888
$ws.Cells.Item(6, 19).Value = 0.05200546038565611
$ws.Cells.Item(6, 20).Value = 0.05200546038565611

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 1.047813333333333
$ws.Cells.Item(7, 8).Value = 3.14344
$ws.Cells.Item(7, 9).Value = 0.05889789011308234
$ws.Cells.Item(7, 10).Value = 0.05889789011308236
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 1.426936666666667
$ws.Cells.Item(7, 14).Value = 4.28081
$ws.Cells.Item(7, 15).Value = 0.04196962907162197
$ws.Cells.Item(7, 16).Value = 0.04196962907162197
$ws.Cells.Item(7, 17).Value = 1.495163265155555
$ws.Cells.Item(7, 18).Value = 13.4564693864
$ws.Cells.Item(7, 19).Value = 0.002471922601147217
$ws.Cells.Item(7, 20).Value = 0.002471922601147218

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 1.047813333333333
$ws.Cells.Item(8, 8).Value = 3.14344
$ws.Cells.Item(8, 9).Value = 0.05889789011308234
$ws.Cells.Item(8, 10).Value = 0.05889789011308236
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 1.067853
$ws.Cells.Item(8, 14).Value = 3.203559
$ws.Cells.Item(8, 15).Value = 0.03140811737476231
$ws.Cells.Item(8, 16).Value = 0.0314081173747623
$ws.Cells.Item(8, 17).Value = 1.11891061144
$ws.Cells.Item(8, 18).Value = 10.07019550296
$ws.Cells.Item(8, 19).Value = 0.001849871845797543
$ws.Cells.Item(8, 20).Value = 0.001849871845797543

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 1.047813333333333
$ws.Cells.Item(9, 8).Value = 3.14344
$ws.Cells.Item(9, 9).Value = 0.05889789011308234
$ws.Cells.Item(9, 10).Value = 0.05889789011308236
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 1.483919333333333
$ws.Cells.Item(9, 14).Value = 4.451758
$ws.Cells.Item(9, 15).Value = 0.04364562593916237
$ws.Cells.Item(9, 16).Value = 0.04364562593916237
$ws.Cells.Item(9, 17).Value = 1.554870463057778
$ws.Cells.Item(9, 18).Value = 13.99383416752
$ws.Cells.Item(9, 19).Value = 0.002570635280481482
$ws.Cells.Item(9, 20).Value = 0.002570635280481483

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 1.312552
$ws.Cells.Item(10, 8).Value = 3.937656
$ws.Cells.Item(10, 9).Value = 0.07377892703252469
$ws.Cells.Item(10, 10).Value = 0.0737789270325247
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 30.020559
$ws.Cells.Item(10, 14).Value = 90.061677
$ws.Cells.Item(10, 15).Value = 0.8829766276144534
$ws.Cells.Item(10, 16).Value = 0.8829766276144534
$ws.Cells.Item(10, 17).Value = 39.403544756568
$ws.Cells.Item(10, 18).Value = 354.631902809112
$ws.Cells.Item(10, 19).Value = 0.06514506818019147
$ws.Cells.Item(10, 20).Value = 0.06514506818019149

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 1.312552
$ws.Cells.Item(11, 8).Value = 3.937656
$ws.Cells.Item(11, 9).Value = 0.07377892703252469
$ws.Cells.Item(11, 10).Value = 0.0737789270325247
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 1.426936666666667
$ws.Cells.Item(11, 14).Value = 4.28081
$ws.Cells.Item(11, 15).Value = 0.04196962907162197
$ws.Cells.Item(11, 16).Value = 0.04196962907162197
$ws.Cells.Item(11, 17).Value = 1.872928575706666
$ws.Cells.Item(11, 18).Value = 16.85635718136
$ws.Cells.Item(11, 19).Value = 0.003096474200857324
$ws.Cells.Item(11, 20).Value = 0.003096474200857325

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 1.312552
$ws.Cells.Item(12, 8).Value = 3.937656
$ws.Cells.Item(12, 9).Value = 0.07377892703252469
$ws.Cells.Item(12, 10).Value = 0.0737789270325247
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 1.067853
$ws.Cells.Item(12, 14).Value = 3.203559
$ws.Cells.Item(12, 15).Value = 0.03140811737476231
$ws.Cells.Item(12, 16).Value = 0.0314081173747623
$ws.Cells.Item(12, 17).Value = 1.401612590856
$ws.Cells.Item(12, 18).Value = 12.614513317704
$ws.Cells.Item(12, 19).Value = 0.002317257200021559
$ws.Cells.Item(12, 20).Value = 0.002317257200021559

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 1.312552
$ws.Cells.Item(13, 8).Value = 3.937656
$ws.Cells.Item(13, 9).Value = 0.07377892703252469
$ws.Cells.Item(13, 10).Value = 0.0737789270325247
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 1.483919333333333
$ws.Cells.Item(13, 14).Value = 4.451758
$ws.Cells.Item(13, 15).Value = 0.04364562593916237
$ws.Cells.Item(13, 16).Value = 0.04364562593916237
$ws.Cells.Item(13, 17).Value = 1.947721288805333
$ws.Cells.Item(13, 18).Value = 17.529491599248
$ws.Cells.Item(13, 19).Value = 0.003220127451454328
$ws.Cells.Item(13, 20).Value = 0.003220127451454328

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 2.458986333333333
$ws.Cells.Item(14, 8).Value = 7.376958999999999
$ws.Cells.Item(14, 9).Value = 0.1382203320409214
$ws.Cells.Item(14, 10).Value = 0.1382203320409214
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 30.020559
$ws.Cells.Item(14, 14).Value = 90.061677
$ws.Cells.Item(14, 15).Value = 0.8829766276144534
$ws.Cells.Item(14, 16).Value = 0.8829766276144534
$ws.Cells.Item(14, 17).Value = 73.82014430002701
$ws.Cells.Item(14, 18).Value = 664.381298700243
$ws.Cells.Item(14, 19).Value = 0.1220453226532427
$ws.Cells.Item(14, 20).Value = 0.1220453226532427

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 2.458986333333333
$ws.Cells.Item(15, 8).Value = 7.376958999999999
$ws.Cells.Item(15, 9).Value = 0.1382203320409214
$ws.Cells.Item(15, 10).Value = 0.1382203320409214
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 1.426936666666667
$ws.Cells.Item(15, 14).Value = 4.28081
$ws.Cells.Item(15, 15).Value = 0.04196962907162197
$ws.Cells.Item(15, 16).Value = 0.04196962907162197
$ws.Cells.Item(15, 17).Value = 3.508817761865555
$ws.Cells.Item(15, 18).Value = 31.57935985679
$ws.Cells.Item(15, 19).Value = 0.005801056065913894
$ws.Cells.Item(15, 20).Value = 0.005801056065913896

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 2.458986333333333
$ws.Cells.Item(16, 8).Value = 7.376958999999999
$ws.Cells.Item(16, 9).Value = 0.1382203320409214
$ws.Cells.Item(16, 10).Value = 0.1382203320409214
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 1.067853
$ws.Cells.Item(16, 14).Value = 3.203559
$ws.Cells.Item(16, 15).Value = 0.03140811737476231
$ws.Cells.Item(16, 16).Value = 0.0314081173747623
$ws.Cells.Item(16, 17).Value = 2.625835933009
$ws.Cells.Item(16, 18).Value = 23.632523397081
$ws.Cells.Item(16, 19).Value = 0.004341240412319878
$ws.Cells.Item(16, 20).Value = 0.004341240412319878

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 2.458986333333333
$ws.Cells.Item(17, 8).Value = 7.376958999999999
$ws.Cells.Item(17, 9).Value = 0.1382203320409214
$ws.Cells.Item(17, 10).Value = 0.1382203320409214
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 1.483919333333333
$ws.Cells.Item(17, 14).Value = 4.451758
$ws.Cells.Item(17, 15).Value = 0.04364562593916237
$ws.Cells.Item(17, 16).Value = 0.04364562593916237
$ws.Cells.Item(17, 17).Value = 3.648937360435778
$ws.Cells.Item(17, 18).Value = 32.840436243922
$ws.Cells.Item(17, 19).Value = 0.006032712909444873
$ws.Cells.Item(17, 20).Value = 0.006032712909444874
